# Regenerate merged AHB files
#
# The header row (row 1) labels the "before" message version columns with a
# "_old" suffix and the "after" message version columns with a "_new" suffix.
# This pass relabels them with the concrete format versions being compared:
#   *_old -> *_FV2210
#   *_new -> *_FV2304
# It also freezes the header row and wraps the sheet's used range in a
# structured table ("Table1") so the headers/filters stay anchored while
# scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:J => "_old" => "_FV2210"
$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"

# Column K ("diff") is unchanged.

# Columns L:U => "_new" => "_FV2304"
$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# Freeze panes above row 2, so the header row stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null

# Wrap the full used range (A1:U54) in a structured table named "Table1",
# matching the header labels above.
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U54"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
